# ---------------------------------------------------------------
# Commit: "Fixed update to excel issue"
#
# 1. Renames the "Requested quantity" header on the "Weekly Quantity"
#    and "Monthly Trend" sheets to "Weekly_PO_Qty" / "Monthly_PO_Qty".
# 2. Adds a new "PO Forecast" sheet (after "Monthly Trend") containing
#    the ds / PO_Forecast / yhat_lower / yhat_upper forecast table.
# ---------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Requested quantity" headers -----------------
$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet, right after "Monthly Trend" ---
$afterSheet  = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast  = $wb.Worksheets.Add($null, $afterSheet)
$wsForecast.Name = "PO Forecast"

# --- 3. Header row values ---------------------------------------
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- 4. Data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ------
$csv = @"
44934.99999999999,263,-458.1347355552525,1089.231699532233
44941.99999999999,266,-513.2155700837236,1007.09629341738
44948.99999999999,269,-479.6675998960729,1017.055891425664
44955.99999999999,273,-535.9929721018798,1026.175923574616
44962.99999999999,276,-494.6511195302305,1014.816544967898
44969.99999999999,280,-418.0941281240117,1046.1512011871
44976.99999999999,283,-466.3644252214917,1026.876085202209
44983.99999999999,286,-417.8293612239113,1090.745970110345
44990.99999999999,290,-487.2561475736513,1064.208670715903
44997.99999999999,293,-477.5945431288764,1092.697467079862
45004.99999999999,297,-520.1575371422132,1028.833801649707
45025.99999999999,307,-440.5174053676748,1017.360873686124
45039.99999999999,314,-466.1882158543564,1088.7400129482
45046.99999999999,317,-445.50908208168,1059.660483621366
45053.99999999999,320,-439.7622936075495,1105.057096863008
45060.99999999999,324,-438.6570231920147,1045.505923769773
45081.99999999999,334,-428.7637908618638,1062.148077868859
45088.99999999999,337,-450.7858446605733,1108.078043433729
45095.99999999999,341,-356.6545628100388,1095.358091952454
45102.99999999999,344,-365.5854241192513,1064.262235427188
45109.99999999999,348,-423.6639607151068,1072.101263252135
45116.99999999999,351,-400.6912199839315,1101.356557967008
45123.99999999999,354,-408.2960869085408,1091.801673384948
45130.99999999999,358,-399.321530953591,1110.866870439488
45137.99999999999,361,-407.3765604853559,1152.139663327238
45144.99999999999,365,-349.7683665415207,1099.298413071868
45151.99999999999,368,-422.5835984945348,1109.770671328268
45158.99999999999,372,-360.5896253836976,1112.782626682152
45165.99999999999,375,-368.7353164801556,1059.75524999599
45172.99999999999,378,-360.4808941326888,1061.100388947686
45179.99999999999,382,-424.4667924635554,1122.978096868642
45186.99999999999,385,-352.2095929839058,1122.086397330909
45193.99999999999,389,-357.0437653008295,1205.14566185536
45200.99999999999,392,-347.2427755611213,1147.90285043484
45207.99999999999,395,-338.5536590971752,1186.854152661453
45214.99999999999,399,-360.5734782700902,1185.278939801267
45221.99999999999,402,-317.4264757633863,1115.984100455499
45228.99999999999,406,-258.5301210350191,1164.236392781532
45235.99999999999,409,-292.167192587026,1146.503006175939
45242.99999999999,412,-326.8057359245191,1163.838011762122
45249.99999999999,416,-326.1870045686604,1096.831864497545
45256.99999999999,419,-328.2188471849653,1138.914470027401
45263.99999999999,423,-360.1422798396339,1243.54190551073
45298.99999999999,440,-252.2715526795573,1207.65504350827
45305.99999999999,443,-335.6428949989225,1205.372684788279
45312.99999999999,446,-340.4124358872751,1206.910722842681
45319.99999999999,450,-371.1665611638832,1199.029490035093
45326.99999999999,453,-268.7390894211337,1212.132989789968
45333.99999999999,457,-268.2513403476991,1219.444386820258
45340.99999999999,460,-291.8211987870018,1243.298504606847
45347.99999999999,464,-289.4590316654548,1184.479893807791
45354.99999999999,467,-320.0354931378259,1180.859258292411
45361.99999999999,470,-246.4648147455729,1201.49974578186
45368.99999999999,474,-241.2054085884555,1228.537924136665
45375.99999999999,477,-305.4101614522805,1222.045358627076
45382.99999999999,481,-264.7684613570395,1235.066939511968
45389.99999999999,484,-265.7102652865966,1252.722539324086
45396.99999999999,487,-194.3334242286177,1265.601625680301
45403.99999999999,491,-234.8746904833094,1260.170458054858
45410.99999999999,494,-247.0709134212381,1212.066725588725
45417.99999999999,498,-200.1495455543614,1259.437727379156
45424.99999999999,501,-240.6497085263271,1241.322431198642
45431.99999999999,504,-211.1951690632821,1260.443518135919
45438.99999999999,508,-274.682019742934,1267.199457018741
45445.99999999999,511,-217.8722988524784,1240.435564577901
45459.99999999999,518,-250.2325773667864,1275.532723780425
45473.99999999999,525,-230.6832101328393,1232.005018459037
45480.99999999999,528,-154.7907394433117,1298.463515990353
45557.99999999999,566,-162.9988038749457,1347.081140237627
45564.99999999999,569,-132.0666996284691,1331.145952281339
45571.99999999999,573,-198.4316594559913,1338.209154719839
45578.99999999999,576,-189.9189486112802,1314.156274108156
45585.99999999999,579,-110.8955688718591,1329.120407805685
45592.99999999999,583,-102.9696743912516,1314.616298934253
45599.99999999999,586,-175.5382394316737,1350.744931395836
45606.99999999999,590,-121.7065135457445,1316.751088702355
45613.99999999999,593,-143.8102479112137,1347.062876468666
45627.99999999999,600,-178.8649458088437,1368.671593240526
45634.99999999999,603,-175.8754735580555,1381.828646245852
45641.99999999999,607,-147.0848884765155,1320.941032515149
45648.99999999999,610,-178.4772438085932,1353.034202577861
45655.99999999999,613,-96.58187227440364,1344.793699503865
45662.99999999999,617,-161.6040030306964,1378.846041818705
45669.99999999999,620,-135.0575982558954,1357.54283741765
45676.99999999999,624,-129.5807642079587,1392.408177388481
45683.99999999999,627,-59.80836841489757,1352.785378903922
"@

$lines = $csv -split "`n"
$n = $lines.Count
$d = New-Object 'object[,]' $n,4
for ($i = 0; $i -lt $n; $i++) {
    $parts = $lines[$i].Trim() -split ","
    $d[$i,0] = [double]$parts[0]
    $d[$i,1] = [double]$parts[1]
    $d[$i,2] = [double]$parts[2]
    $d[$i,3] = [double]$parts[3]
}

$lastRow = 1 + $n
$dataRange = $wsForecast.Range("A2:D" + $lastRow)
$dataRange.Value = $d

# --- 5. Re-use the exact formatting already used elsewhere in the
#        workbook (bold/centered header style, and the date
#        number-format used for the "ds"/date columns) instead of
#        inventing new style records. ---------------------------
$headerFormatSource = $wsWeekly.Range("A1:B1")
$headerFormatSource.Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$dateFormatSource = $wsWeekly.Range("A2:A79")
$dateFormatSource.Copy()
$wsForecast.Range("A2:A" + $lastRow).PasteSpecial(-4122)

$excel.CutCopyMode = 0
